$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap RenderToken / VeChain rows (48 <-> 49) for Coin name and Link ---
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"

# --- Price (column D) updates: force text to avoid numeric auto-conversion, then restore default style ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.047.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.450.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "483.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.455.08"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0965"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.859.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.127.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.454.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "313.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "57.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.408"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.560.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "147.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.855"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.75"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.997"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0551"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.600"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0921"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "257.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0227"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.848.62"
$ws.Range("D51").Style = "Normal"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("E3").Value = "  -2.12%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("E6").Value = "  +4.91%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -1.12%  "
$ws.Range("E9").Value = "  -2.78%  "
$ws.Range("E10").Value = "  +5.14%  "
$ws.Range("E11").Value = "  -2.53%  "
$ws.Range("E12").Value = "  +1.60%  "
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("E14").Value = "  -3.31%  "
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("E16").Value = "  +2.20%  "
$ws.Range("E17").Value = "  -3.51%  "
$ws.Range("E18").Value = "  -2.39%  "
$ws.Range("E19").Value = "  +3.82%  "
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("E21").Value = "  -3.12%  "
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("E23").Value = "  +3.01%  "
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("E25").Value = "  +1.79%  "
$ws.Range("E26").Value = "  -0.89%  "
$ws.Range("E27").Value = "  -2.55%  "
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("E29").Value = "  +1.81%  "
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("E32").Value = "  -1.68%  "
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("E35").Value = "  -2.26%  "
$ws.Range("E36").Value = "  +4.38%  "
$ws.Range("E37").Value = "  -0.78%  "
$ws.Range("E38").Value = "  +1.66%  "
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("E40").Value = "  +2.83%  "
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("E43").Value = "  -3.10%  "
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("E45").Value = "  +1.66%  "
$ws.Range("E46").Value = "  +5.40%  "
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("E48").Value = "  +1.99%  "
$ws.Range("E49").Value = "  +4.46%  "
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("E51").Value = "  -6.69%  "
